$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.74
$ws.Range("B8").Value = 6.374000000000001
$ws.Range("B10").Value = 6.096
$ws.Range("B12").Value = 5.258
$ws.Range("C12").Value = -10.95
$ws.Range("C15").Value = -13.364
$ws.Range("C17").Value = -13.267
$ws.Range("B18").Value = 5.583
$ws.Range("C26").Value = -13.15
$ws.Range("C27").Value = -13.416
$ws.Range("C28").Value = -12.77
$ws.Range("B37").Value = 8.190999999999999
$ws.Range("C37").Value = -12.055
$ws.Range("C47").Value = -12.542
$ws.Range("B55").Value = 4.865
$ws.Range("C65").Value = -12.161
$ws.Range("B68").Value = 4.746
$ws.Range("C73").Value = -12.215
$ws.Range("B77").Value = 6.209999999999999
$ws.Range("B78").Value = 7.708
$ws.Range("B81").Value = 5.931
$ws.Range("B82").Value = 5.554
$ws.Range("C84").Value = -13.645
$ws.Range("C85").Value = -12.5
$ws.Range("C93").Value = -10.438
$ws.Range("C95").Value = -12.451
$ws.Range("C98").Value = -13.276
$ws.Range("C99").Value = -11.696
$ws.Range("C101").Value = -12.613
